# EIA Table 4.8.A monthly refresh: October 2016/2015 -> November 2016/2015
# (EPM_2016_11P run). Updates the report title, the two period-header
# columns, and the state/region figures that moved between the October
# and November data pulls.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Helper: write literal text into a cell without letting Excel's
# type-inference reinterpret month/year-looking text (e.g. "November 2016")
# as a date serial. Building it as a formula that returns the string, then
# collapsing the formula to its cached value via Copy/PasteSpecial-Values,
# commits the text as a plain string without flipping the cell's number
# format or style.
function Set-LiteralText {
    param($rng, [string]$text)
    $rng.Formula = "=""" + $text.Replace('"', '""') + """"
    $rng.Copy() | Out-Null
    $rng.PasteSpecial(-4163) | Out-Null
}

# --- Title (row 1) ---
$ws.Range("A1").Value = "Table 4.8.A. Receipts of Petroleum Coke Delivered for Electricity Generation by State, November 2016 and 2015"

# --- Column period headers (row 5) ---
# B/E/G/I/K = current-month column ("...2016"); C/F/H/J/L = prior-year column ("...2015")
Set-LiteralText $ws.Range("B5") "November 2016"
Set-LiteralText $ws.Range("C5") "November 2015"
Set-LiteralText $ws.Range("E5") "November 2016"
Set-LiteralText $ws.Range("F5") "November 2015"
Set-LiteralText $ws.Range("G5") "November 2016"
Set-LiteralText $ws.Range("H5") "November 2015"
Set-LiteralText $ws.Range("I5") "November 2016"
Set-LiteralText $ws.Range("J5") "November 2015"
Set-LiteralText $ws.Range("K5") "November 2016"
Set-LiteralText $ws.Range("L5") "November 2015"

# --- Data refresh ---

# Row 17: East North Central
$ws.Range("B17").Value = 73
$ws.Range("C17").Value = 143
$ws.Range("D17").Value = -0.49
$ws.Range("E17").Value = 27
$ws.Range("F17").Value = 78
$ws.Range("G17").Value = 46
$ws.Range("H17").Value = 59
$ws.Range("L17").Value = 5

# Row 19: Indiana
$ws.Range("C19").Value = 56
$ws.Range("F19").Value = 56

# Row 20: Michigan
$ws.Range("B20").Value = 19
$ws.Range("C20").Value = 19
$ws.Range("D20").Value = -0.002
$ws.Range("E20").Value = 19
$ws.Range("F20").Value = 19

# Row 21: Ohio
$ws.Range("B21").Value = 46
$ws.Range("C21").Value = 59
$ws.Range("D21").Value = -0.22
$ws.Range("G21").Value = 46
$ws.Range("H21").Value = 59

# Row 22: Wisconsin
$ws.Range("B22").Value = 8
$ws.Range("C22").Value = 9
$ws.Range("D22").Value = -0.092
$ws.Range("E22").Value = 8
$ws.Range("F22").Value = 4
$ws.Range("L22").Value = 5

# Row 31: South Atlantic
$ws.Range("B31").Value = 84
$ws.Range("C31").Value = 96
$ws.Range("D31").Value = -0.12
$ws.Range("E31").Value = 76
$ws.Range("F31").Value = 96

# Row 34: Florida
$ws.Range("B34").Value = 76
$ws.Range("C34").Value = 96
$ws.Range("D34").Value = -0.21
$ws.Range("E34").Value = 76
$ws.Range("F34").Value = 96

# Row 41: East South Central
$ws.Range("B41").Value = 13
$ws.Range("C41").Value = 48
$ws.Range("D41").Value = -0.73
$ws.Range("E41").Value = 13
$ws.Range("F41").Value = 48

# Row 43: Kentucky
$ws.Range("B43").Value = 13
$ws.Range("C43").Value = 48
$ws.Range("D43").Value = -0.73
$ws.Range("E43").Value = 13
$ws.Range("F43").Value = 48

# Row 46: West South Central
$ws.Range("B46").Value = 163
$ws.Range("C46").Value = 132
$ws.Range("D46").Value = 0.23
$ws.Range("E46").Value = 163
$ws.Range("F46").Value = 132

# Row 48: Louisiana
$ws.Range("B48").Value = 163
$ws.Range("C48").Value = 132
$ws.Range("D48").Value = 0.23
$ws.Range("E48").Value = 163
$ws.Range("F48").Value = 132

# Row 67: U.S. Total
$ws.Range("B67").Value = 333
$ws.Range("C67").Value = 429
$ws.Range("D67").Value = -0.22
$ws.Range("E67").Value = 279
$ws.Range("F67").Value = 354
$ws.Range("G67").Value = 46
$ws.Range("H67").Value = 59
$ws.Range("L67").Value = 15
